$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 2
$ws.Range("H2").Value = 1587416.1
$ws.Range("I2").Value = 2777865.2
$ws.Range("J2").Value = 150.66667
$ws.Range("K2").Value = 2777865.2
$ws.Range("L2").Value = 150.66667
$ws.Range("M2").Value = -2777752.2
$ws.Range("N2").Value = -376.66667

# ALC row 21
$ws.Range("H21").Value = 40004.2
$ws.Range("J21").Value = 40004.2
$ws.Range("L21").Value = 40004.2
$ws.Range("N21").Value = -40940.2

# ALC row 23
$ws.Range("H23").Value = 40004.2
$ws.Range("J23").Value = 40004.2
$ws.Range("L23").Value = 40004.2
$ws.Range("N23").Value = -40472.2

# ALC row 29
$ws.Range("H29").Value = 1500
$ws.Range("I29").Value = 1500
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 4500
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -4219
$ws.Range("N29").ClearContents()

# ALC row 38
$ws.Range("H38").Value = 78.86667
$ws.Range("I38").Value = 78.86667
$ws.Range("K38").Value = 236.60001
$ws.Range("M38").Value = 135.39999

# ALC row 112
$ws.Range("H112").Value = 1059.3617
$ws.Range("J112").Value = 1079.3182
$ws.Range("L112").Value = 3237.9546
$ws.Range("N112").Value = -5453.9546

# ALC row 137
$ws.Range("H137").Value = 17242570
$ws.Range("I137").Value = 26316594
$ws.Range("J137").Value = 1920.6
$ws.Range("K137").Value = 78949782
$ws.Range("L137").Value = 5761.799999999999
$ws.Range("M137").Value = -78947232
$ws.Range("N137").Value = -10861.8

$ws = $wb.Worksheets.Item("ARM")
# ARM row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

# ARM row 9
$ws.Range("H9").Value = 8220
$ws.Range("J9").Value = 8220
$ws.Range("L9").Value = 8220
$ws.Range("N9").Value = -8560

# ARM row 20
$ws.Range("H20").Value = 8220
$ws.Range("J20").Value = 8220
$ws.Range("L20").Value = 8220
$ws.Range("N20").Value = -8760

# ARM row 23
$ws.Range("H23").Value = 12001.2
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 12001.2
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 12001.2
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -12519.2

# ARM row 32
$ws.Range("H32").Value = 14558.789
$ws.Range("I32").Value = 11179.713
$ws.Range("J32").Value = 41591.4
$ws.Range("K32").Value = 11179.713
$ws.Range("L32").Value = 41591.4
$ws.Range("M32").Value = -10892.713
$ws.Range("N32").Value = -42165.4

# ARM row 37
$ws.Range("H37").Value = 10472.5
$ws.Range("I37").Value = 4611.6665
$ws.Range("J37").Value = 16333.333
$ws.Range("K37").Value = 4611.6665
$ws.Range("L37").Value = 16333.333
$ws.Range("M37").Value = -4338.6665
$ws.Range("N37").Value = -16879.333

# ARM row 44
$ws.Range("H44").Value = 25000
$ws.Range("J44").Value = 25000
$ws.Range("L44").Value = 25000
$ws.Range("N44").Value = -25976

# ARM row 55
$ws.Range("H55").Value = 25000
$ws.Range("J55").Value = 25000
$ws.Range("L55").Value = 25000
$ws.Range("N55").Value = -25630

# ARM row 61
$ws.Range("H61").Value = 1815.2106
$ws.Range("I61").Value = 1529.9592
$ws.Range("K61").Value = 1529.9592
$ws.Range("M61").Value = -1317.9592

# ARM row 136
$ws.Range("H136").Value = 1815.2106
$ws.Range("I136").Value = 1529.9592
$ws.Range("K136").Value = 4589.8776
$ws.Range("M136").Value = -2039.8776

$ws = $wb.Worksheets.Item("BSM")
# BSM row 105
$ws.Range("H105").Value = 1264436.5
$ws.Range("I105").Value = 1624939.8
$ws.Range("K105").Value = 1624939.8
$ws.Range("M105").Value = -1623192.8

# BSM row 107
$ws.Range("H107").Value = 1108
$ws.Range("I107").Value = 1065.75
$ws.Range("J107").Value = 1333.3334
$ws.Range("K107").Value = 1065.75
$ws.Range("L107").Value = 1333.3334
$ws.Range("M107").Value = 854.25
$ws.Range("N107").Value = -5173.3334

$ws = $wb.Worksheets.Item("CRP")
# CRP row 22
$ws.Range("H22").Value = 693.3333
$ws.Range("I22").Value = 441.375
$ws.Range("J22").Value = 894.9
$ws.Range("K22").Value = 441.375
$ws.Range("L22").Value = 894.9
$ws.Range("M22").Value = -91.375
$ws.Range("N22").Value = -1594.9

# CRP row 31
$ws.Range("H31").Value = 13892311
$ws.Range("I31").Value = 27779046
$ws.Range("J31").Value = 5576.528
$ws.Range("K31").Value = 27779046
$ws.Range("L31").Value = 5576.528
$ws.Range("M31").Value = -27778751
$ws.Range("N31").Value = -6166.528

# CRP row 34
$ws.Range("H34").Value = 13892311
$ws.Range("I34").Value = 27779046
$ws.Range("J34").Value = 5576.528
$ws.Range("K34").Value = 27779046
$ws.Range("L34").Value = 5576.528
$ws.Range("M34").Value = -27778844
$ws.Range("N34").Value = -5980.528

# CRP row 107
$ws.Range("H107").Value = 1174.44
$ws.Range("I107").Value = 242.1
$ws.Range("J107").Value = 1796
$ws.Range("K107").Value = 242.1
$ws.Range("L107").Value = 1796
$ws.Range("M107").Value = 1677.9
$ws.Range("N107").Value = -5636

# CRP row 138
$ws.Range("H138").Value = 64945
$ws.Range("J138").Value = 65651.42999999999
$ws.Range("L138").Value = 65651.42999999999
$ws.Range("N138").Value = -75931.42999999999

# CRP row 139
$ws.Range("H139").Value = 53300
$ws.Range("J139").Value = 49950
$ws.Range("L139").Value = 49950
$ws.Range("N139").Value = -60230

# CRP row 140
$ws.Range("H140").Value = 50066.5
$ws.Range("J140").Value = 49938
$ws.Range("L140").Value = 49938
$ws.Range("N140").Value = -60298

# CRP row 141
$ws.Range("H141").Value = 29464.166
$ws.Range("J141").Value = 29464.166
$ws.Range("L141").Value = 29464.166
$ws.Range("N141").Value = -39824.166

$ws = $wb.Worksheets.Item("CUL")
# CUL row 92
$ws.Range("H92").Value = 506.08334
$ws.Range("I92").Value = 638
$ws.Range("J92").Value = 411.85715
$ws.Range("K92").Value = 1914
$ws.Range("L92").Value = 1235.57145
$ws.Range("M92").Value = -666
$ws.Range("N92").Value = -3731.57145

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws.Range("H80").Value = 75371.5
$ws.Range("I80").Value = 2639
$ws.Range("J80").Value = 115778.445
$ws.Range("K80").Value = 2639
$ws.Range("L80").Value = 115778.445
$ws.Range("M80").Value = -1641
$ws.Range("N80").Value = -117774.445

# GSM row 83
$ws.Range("H83").Value = 75371.5
$ws.Range("I83").Value = 2639
$ws.Range("J83").Value = 115778.445
$ws.Range("K83").Value = 13195
$ws.Range("L83").Value = 578892.2250000001
$ws.Range("M83").Value = -8203
$ws.Range("N83").Value = -588876.2250000001

# GSM row 126
$ws.Range("H126").Value = 20835724
$ws.Range("I126").Value = 27780000
$ws.Range("J126").Value = 2900
$ws.Range("K126").Value = 83340000
$ws.Range("L126").Value = 8700
$ws.Range("M126").Value = -83337530
$ws.Range("N126").Value = -13640

$ws = $wb.Worksheets.Item("LTW")
# LTW row 9
$ws.Range("H9").Value = 12025.1
$ws.Range("I9").Value = 166.88889
$ws.Range("J9").Value = 21727.273
$ws.Range("K9").Value = 166.88889
$ws.Range("L9").Value = 21727.273
$ws.Range("M9").Value = 57.11111
$ws.Range("N9").Value = -22175.273

# LTW row 46
$ws.Range("H46").Value = 1026.2
$ws.Range("I46").Value = 792.5
$ws.Range("J46").Value = 1182
$ws.Range("K46").Value = 792.5
$ws.Range("L46").Value = 1182
$ws.Range("M46").Value = -604.5
$ws.Range("N46").Value = -1558

# LTW row 82
$ws.Range("H82").Value = 3250
$ws.Range("I82").Value = 3250
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3250
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2889
$ws.Range("N82").ClearContents()

# LTW row 85
$ws.Range("H85").Value = 3250
$ws.Range("I85").Value = 3250
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3250
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2002
$ws.Range("N85").ClearContents()
